$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
Write-Host $ws.Name
